$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 164, shifting existing rows
# (164-215) down to (165-216).
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new price record.
$ws.Range("A164").Value = 11
$ws.Range("B164").Value = 'Vega Monumental Concepción'
$ws.Range("C164").Value = 'Bíobío'
$ws.Range("D164").Value = 44588
$ws.Range("E164").Value = 8
$ws.Range("F164").Value = 'Fruta'
$ws.Range("G164").Value = 100102
$ws.Range("H164").Value = 'Cítricos'
$ws.Range("I164").Value = 100102005
$ws.Range("J164").Value = 'Naranja'
$ws.Range("K164").Value = 'Valencia'
$ws.Range("L164").Value = 'Primera'
$ws.Range("M164").Value = 220
$ws.Range("N164").Value = 9000
$ws.Range("O164").Value = 9500
$ws.Range("P164").Value = 9227
$ws.Range("Q164").Value = '$/caja 18 kilos importada'
$ws.Range("R164").Value = "Región de O'Higgins"
$ws.Range("S164").Value = 513
$ws.Range("T164").Value = 18
